$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.403.41'
$ws.Range("E2").Value = '  -0.97%  '
$ws.Range("D3").Value = '1.564.65'
$ws.Range("E3").Value = '  -1.31%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '''208.95'
$ws.Range("E5").Value = '  +1.10%  '
$ws.Range("E6").Value = '  -0.39%  '
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").Value = '''21.98'
$ws.Range("E8").Value = '  -1.39%  '
$ws.Range("E9").Value = '  -1.85%  '
$ws.Range("E10").Value = '  +0.02%  '
$ws.Range("D11").Value = '''0.0865'
$ws.Range("E11").Value = '  -0.37%  '
$ws.Range("D12").Value = '1.789.39'
$ws.Range("E12").Value = '  -1.19%  '
$ws.Range("D13").Value = '1.565.22'
$ws.Range("E13").Value = '  -1.29%  '
$ws.Range("D14").Value = '''3.82'
$ws.Range("E14").Value = '  -1.19%  '
$ws.Range("D15").Value = '''0.516'
$ws.Range("E15").Value = '  -2.90%  '
$ws.Range("D16").Value = '''63.48'
$ws.Range("E16").Value = '  +0.60%  '
$ws.Range("D17").Value = '27.392.55'
$ws.Range("E17").Value = '  -0.97%  '
$ws.Range("D18").Value = '''212.77'
$ws.Range("E18").Value = '  -2.78%  '
$ws.Range("E19").Value = '  -0.68%  '
$ws.Range("D20").Value = '''7.25'
$ws.Range("E20").Value = '  -0.90%  '
$ws.Range("E21").Value = '  -0.17%  '
$ws.Range("E22").Value = '  -0.59%  '
$ws.Range("D23").Value = '''9.54'
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("E24").Value = '  +1.40%  '
$ws.Range("D25").Value = '''153.38'
$ws.Range("E25").Value = '  -0.14%  '
$ws.Range("D27").Value = '''6.72'
$ws.Range("E27").Value = '  -0.16%  '
$ws.Range("D28").Value = '''14.96'
$ws.Range("E28").Value = '  -0.85%  '
$ws.Range("E29").Value = '  -2.16%  '
$ws.Range("E30").Value = '  +0.07%  '
$ws.Range("E31").Value = '  +0.95%  '
$ws.Range("E32").Value = '  -0.64%  '
$ws.Range("D33").Value = '1.370.24'
$ws.Range("E33").Value = '  -0.83%  '
$ws.Range("E34").Value = '  +0.15%  '
$ws.Range("E35").Value = '  +1.47%  '
$ws.Range("D36").Value = '''0.964'
$ws.Range("E36").Value = '  -0.32%  '
$ws.Range("E37").Value = '  -0.42%  '
$ws.Range("E38").Value = '  +0.96%  '
$ws.Range("D39").Value = '''0.530'
$ws.Range("E39").Value = '  -2.00%  '
$ws.Range("D40").Value = '''0.822'
$ws.Range("E40").Value = '  +0.10%  '
$ws.Range("E41").Value = '  -0.16%  '
$ws.Range("D42").Value = '''0.973'
$ws.Range("E42").Value = '  -0.65%  '
$ws.Range("E43").Value = '  +0.65%  '
$ws.Range("D44").Value = '''63.95'
$ws.Range("E44").Value = '  +0.12%  '
$ws.Range("E45").Value = '  +0.58%  '
$ws.Range("E46").Value = '  -1.02%  '
$ws.Range("D47").Value = '1.700.68'
$ws.Range("E47").Value = '  -1.16%  '
$ws.Range("D48").Value = '''85.58'
$ws.Range("E48").Value = '  -2.44%  '
$ws.Range("D49").Value = '0.0₇0993'
$ws.Range("E49").Value = '  -2.03%  '
$ws.Range("D50").Value = '''0.0953'
$ws.Range("E50").Value = '  -2.25%  '
$ws.Range("E51").Value = '  -0.90%  '
